$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply per-row updates to reflect refreshed crypto price/volume data.
# Columns B (Coin) and C (Link) only change where rows were reordered;
# D (Price) and E (Volume 1h) are refreshed text values, force text format
# so Excel does not silently reinterpret numeric-looking strings as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.198.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.783.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.98"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.78"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.07%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.041.76"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.01"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.785.70"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.103.66"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.623"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.13"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0785"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.98"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.26%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.73%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.72"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.69%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.62%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.73"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.64"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.444.26"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.73%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.43"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.90%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.05%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.54%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.10"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.98%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.37"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.923"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.29%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.56"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.99%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.16%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.943.38"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.60"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.22%  "
